$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to text format so numeric-looking strings
# (e.g. "317.52", "28.464.85") are stored as text, matching the source data.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '28.464.85'
$ws.Range('E2').Value = '  -0.35%  '
$ws.Range('D3').Value = '1.878.59'
$ws.Range('E3').Value = '  -0.88%  '
$ws.Range('D4').Value = '1.018'
$ws.Range('E4').Value = '  -2.48%  '
$ws.Range('D5').Value = '317.52'
$ws.Range('E5').Value = '  -0.33%  '
$ws.Range('E6').Value = '  -0.24%  '
$ws.Range('D7').Value = '0.5124'
$ws.Range('E7').Value = '  -0.92%  '
$ws.Range('D8').Value = '0.3958'
$ws.Range('E8').Value = '  +0.27%  '
$ws.Range('D9').Value = '0.08449'
$ws.Range('E9').Value = '  +0.46%  '
$ws.Range('D10').Value = '1.112'
$ws.Range('E10').Value = '  -1.61%  '
$ws.Range('B11').Value = 'Polkadot'
$ws.Range('C11').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D11').Value = '6.265'
$ws.Range('E11').Value = '  -0.26%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.870.06'
$ws.Range('E12').Value = '  -2.31%  '
$ws.Range('B13').Value = 'Solana'
$ws.Range('C13').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D13').Value = '20.50'
$ws.Range('E13').Value = '  -0.32%  '
$ws.Range('B14').Value = 'Chainlink'
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D14').Value = '7.253'
$ws.Range('E14').Value = '  -0.78%  '
$ws.Range('B15').Value = 'BinanceUSD'
$ws.Range('C15').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D15').Value = '1.022'
$ws.Range('E15').Value = '  -1.78%  '
$ws.Range('B16').Value = 'ShibaInu'
$ws.Range('C16').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D16').Value = '0.00001111'
$ws.Range('E16').Value = '  -0.05%  '
$ws.Range('B17').Value = 'Litecoin'
$ws.Range('C17').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D17').Value = '91.21'
$ws.Range('E17').Value = '  -0.04%  '
$ws.Range('B18').Value = 'TRON'
$ws.Range('C18').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D18').Value = '0.06795'
$ws.Range('E18').Value = '  +0.41%  '
$ws.Range('B19').Value = 'Avalanche'
$ws.Range('C19').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D19').Value = '17.75'
$ws.Range('E19').Value = '  -1.02%  '
$ws.Range('B20').Value = 'Dai'
$ws.Range('C20').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D20').Value = '1.018'
$ws.Range('E20').Value = '  -0.33%  '
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').Value = '5.959'
$ws.Range('E21').Value = '  -1.96%  '
$ws.Range('B22').Value = 'WrappedBTC'
$ws.Range('C22').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D22').Value = '28.476.06'
$ws.Range('E22').Value = '  -0.10%  '
$ws.Range('B23').Value = 'Cosmos'
$ws.Range('C23').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D23').Value = '11.19'
$ws.Range('E23').Value = '  -0.34%  '
$ws.Range('B24').Value = 'Toncoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D24').Value = '2.283'
$ws.Range('E24').Value = '  -0.95%  '
$ws.Range('B25').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C25').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D25').Value = '2.086.23'
$ws.Range('E25').Value = '  -1.80%  '
$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').Value = '161.73'
$ws.Range('E26').Value = '  -0.15%  '
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').Value = '20.81'
$ws.Range('E27').Value = '  -0.53%  '
$ws.Range('B28').Value = 'LidoDAOToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D28').Value = '2.366'
$ws.Range('E28').Value = '  -3.79%  '
$ws.Range('B29').Value = 'BitcoinCash'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D29').Value = '126.71'
$ws.Range('E29').Value = '  -0.44%  '
$ws.Range('B30').Value = 'Stellar'
$ws.Range('C30').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D30').Value = '0.1057'
$ws.Range('E30').Value = '  -0.80%  '
$ws.Range('B31').Value = 'ImmutableX'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D31').Value = '1.045'
$ws.Range('E31').Value = '  -0.06%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').Value = '5.782'
$ws.Range('E32').Value = '  -2.88%  '
$ws.Range('B33').Value = 'HuobiToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D33').Value = '3.624'
$ws.Range('E33').Value = '  -1.09%  '
$ws.Range('B34').Value = 'VeChain'
$ws.Range('C34').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D34').Value = '0.02434'
$ws.Range('E34').Value = '  -0.98%  '
$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D35').Value = '0.06497'
$ws.Range('E35').Value = '  -1.97%  '
$ws.Range('B36').Value = 'Algorand'
$ws.Range('C36').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D36').Value = '0.2178'
$ws.Range('E36').Value = '  -2.27%  '
$ws.Range('B37').Value = 'FraxShare'
$ws.Range('C37').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D37').Value = '8.913'
$ws.Range('E37').Value = '  -5.80%  '
$ws.Range('B38').Value = 'TrustWalletToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D38').Value = '1.265'
$ws.Range('E38').Value = '  +0.57%  '
$ws.Range('B39').Value = 'ARBITRUM'
$ws.Range('C39').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D39').Value = '1.190'
$ws.Range('E39').Value = '  -0.03%  '
$ws.Range('B40').Value = 'TheSandbox'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D40').Value = '0.6399'
$ws.Range('E40').Value = '  -2.18%  '
$ws.Range('B41').Value = 'InternetComputer(DFINITY)'
$ws.Range('C41').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D41').Value = '5.011'
$ws.Range('E41').Value = '  -0.59%  '
$ws.Range('B42').Value = 'Aptos'
$ws.Range('C42').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D42').Value = '11.24'
$ws.Range('E42').Value = '  -0.01%  '
$ws.Range('B43').Value = 'Frax'
$ws.Range('C43').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D43').Value = '1.017'
$ws.Range('E43').Value = '  +0.00%  '
$ws.Range('B44').Value = 'Decentraland'
$ws.Range('C44').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D44').Value = '0.6058'
$ws.Range('E44').Value = '  -1.53%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').Value = '12.97'
$ws.Range('E45').Value = '  -1.97%  '
$ws.Range('B46').Value = 'PancakeSwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D46').Value = '3.724'
$ws.Range('E46').Value = '  -0.15%  '
$ws.Range('B47').Value = 'NEARProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D47').Value = '1.994'
$ws.Range('E47').Value = '  -1.22%  '
$ws.Range('B48').Value = 'WEMIXTOKEN'
$ws.Range('C48').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D48').Value = '1.206'
$ws.Range('E48').Value = '  -6.93%  '
$ws.Range('B49').Value = 'Quant'
$ws.Range('C49').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D49').Value = '122.27'
$ws.Range('E49').Value = '  +0.20%  '
$ws.Range('B50').Value = 'EOS'
$ws.Range('C50').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D50').Value = '1.207'
$ws.Range('E50').Value = '  -2.77%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').Value = '0.06850'
$ws.Range('E51').Value = '  -1.25%  '

# Remove the explicit text-format style again so cells keep their original
# (unstyled) appearance, now that the values are safely stored as text.
$ws.Range("D2:D51").ClearFormats()
